$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update team labels for rows 3 and 4 (order in shared-string table swaps
#     so the text shown for these two rows trades places) ---
$ws.Range("B3").Value = "Xinghao_s2l"
$ws.Range("B4").Value = "Kurni_s2l"

# --- Update Weekly Pending Total (C) and Repayment (D) figures ---
$ws.Range("C2").Value = 1398025598
$ws.Range("D2").Value = 142152999

$ws.Range("C3").Value = 2053407690
$ws.Range("D3").Value = 208243959

$ws.Range("C4").Value = 5020565466
$ws.Range("D4").Value = 503680243

$ws.Range("C5").Value = 1388402014
$ws.Range("D5").Value = 137532559

$ws.Range("C6").Value = 5122108931
$ws.Range("D6").Value = 431227148

$ws.Range("C7").Value = 2667443011
$ws.Range("D7").Value = 215797676

$ws.Range("C8").Value = 5003636165
$ws.Range("D8").Value = 398981630

# --- Recovery rate (E) is D/C ; ranks (F) stay sequential 1-7, unchanged ---

# --- Autofit columns A:F to best-fit widths (matching Excel's "best fit" sizing) ---
$ws.Columns.Item(1).ColumnWidth = 4.751
$ws.Columns.Item(2).ColumnWidth = 11.084333333333333
$ws.Columns.Item(3).ColumnWidth = 23.584333333333333
$ws.Columns.Item(4).ColumnWidth = 10.251
$ws.Columns.Item(5).ColumnWidth = 12.417666666666666
$ws.Columns.Item(6).ColumnWidth = 4.417666666666667

# --- Move the active selection from J8 to K8 ---
$ws.Range("K8").Select() | Out-Null
